$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82
$ws.Range("K82").Value = 15385.2861
$ws.Range("I82").Value = 5128.4287
$ws.Range("M82").Value = -14979.2861
$ws.Range("H82").Value = 5128.4287

# Row 17
$ws.Range("L17").Value = 3127.3077
$ws.Range("N17").Value = -3463.3077
$ws.Range("H17").Value = 1030.5238
$ws.Range("J17").Value = 1042.4359

# Row 61
$ws.Range("I61").Value = 772.3
$ws.Range("M61").Value = -2144.9
$ws.Range("H61").Value = 772.3
$ws.Range("K61").Value = 2316.9

# Row 85
$ws.Range("H85").Value = 5128.4287
$ws.Range("M85").Value = -13981.2861
$ws.Range("K85").Value = 15385.2861
$ws.Range("I85").Value = 5128.4287

# Row 15
$ws.Range("K15").Value = 6779.400000000001
$ws.Range("H15").Value = 2259.8
$ws.Range("I15").Value = 2259.8
$ws.Range("M15").Value = -6610.400000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("I97").Value = 977.7646999999999
$ws.Range("K97").Value = 977.7646999999999
$ws.Range("H97").Value = 1127.6086
$ws.Range("M97").Value = -481.7646999999999

# Row 136
$ws.Range("H136").Value = 1519218.1
$ws.Range("M136").Value = -6257612.4
$ws.Range("N136").Value = -22733.5005
$ws.Range("J136").Value = 5877.8335
$ws.Range("L136").Value = 17633.5005
$ws.Range("I136").Value = 2086720.8
$ws.Range("K136").Value = 6260162.4

# Row 74
$ws.Range("K74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H74").Value = 9198.5

# Row 110
$ws.Range("I110").Value = 1079.9375
$ws.Range("L110").Value = 600
$ws.Range("J110").Value = 600
$ws.Range("N110").Value = -4690
$ws.Range("H110").Value = 1051.7059
$ws.Range("K110").Value = 1079.9375
$ws.Range("M110").Value = 965.0625

# Row 122
$ws.Range("M122").Value = -4184.3125
$ws.Range("H122").Value = 2162.7
$ws.Range("J122").Value = 1967.75
$ws.Range("K122").Value = 6634.3125
$ws.Range("L122").Value = 5903.25
$ws.Range("N122").Value = -10803.25
$ws.Range("I122").Value = 2211.4375

# Row 77
$ws.Range("M77").ClearContents()
$ws.Range("K77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("H77").Value = 9198.5

# Row 61
$ws.Range("J61").Value = 5877.8335
$ws.Range("L61").Value = 5877.8335
$ws.Range("I61").Value = 2086720.8
$ws.Range("M61").Value = -2086508.8
$ws.Range("H61").Value = 1519218.1
$ws.Range("N61").Value = -6301.8335
$ws.Range("K61").Value = 2086720.8

# Row 117
$ws.Range("L117").Value = 79999
$ws.Range("H117").Value = 79999
$ws.Range("J117").Value = 79999
$ws.Range("N117").Value = -89177

# Row 102
$ws.Range("H102").Value = 22794.264
$ws.Range("I102").Value = 26130.875
$ws.Range("K102").Value = 26130.875
$ws.Range("M102").Value = -24508.875

# Row 132
$ws.Range("M132").Value = -2375639.3
$ws.Range("K132").Value = 2378169.3
$ws.Range("H132").Value = 643705.0600000001
$ws.Range("I132").Value = 792723.1

# Row 32
$ws.Range("L32").Value = 98.5
$ws.Range("I32").Value = 5467147
$ws.Range("N32").Value = -672.5
$ws.Range("K32").Value = 5467147
$ws.Range("H32").Value = 5293590
$ws.Range("M32").Value = -5466860
$ws.Range("J32").Value = 98.5

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("K82").Value = 63000
$ws.Range("I82").Value = 63000
$ws.Range("M82").Value = -62617
$ws.Range("L82").Value = 79213
$ws.Range("J82").Value = 79213
$ws.Range("N82").Value = -79979
$ws.Range("H82").Value = 76896.86

# Row 134
$ws.Range("H134").Value = 1373658
$ws.Range("M134").Value = -24155334
$ws.Range("J134").Value = 260497.22
$ws.Range("L134").Value = 781491.66
$ws.Range("N134").Value = -786561.66
$ws.Range("I134").Value = 8052623
$ws.Range("K134").Value = 24157869

# Row 100
$ws.Range("J100").Value = 34500
$ws.Range("N100").Value = -36664
$ws.Range("L100").Value = 34500
$ws.Range("H100").Value = 34500

# Row 20
$ws.Range("H20").Value = 1976.9166
$ws.Range("I20").Value = 1262.6666
$ws.Range("K20").Value = 1262.6666
$ws.Range("M20").Value = -1015.6666

# Row 85
$ws.Range("N85").Value = -81865
$ws.Range("H85").Value = 76896.86
$ws.Range("L85").Value = 79213
$ws.Range("M85").Value = -61674
$ws.Range("K85").Value = 63000
$ws.Range("J85").Value = 79213
$ws.Range("I85").Value = 63000

# Row 105
$ws.Range("H105").Value = 1749.25
$ws.Range("I105").Value = 1749.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1749.25
$ws.Range("M105").Value = -2.25
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 120
$ws.Range("N120").Value = -63899.832
$ws.Range("J120").Value = 56641.832
$ws.Range("H120").Value = 56641.832
$ws.Range("L120").Value = 56641.832

# Row 31
$ws.Range("H31").Value = 108411.39
$ws.Range("K31").Value = 224465.22
$ws.Range("M31").Value = -224170.22
$ws.Range("J31").Value = 22898.053
$ws.Range("N31").Value = -23488.053
$ws.Range("L31").Value = 22898.053
$ws.Range("I31").Value = 224465.22

# Row 122
$ws.Range("M122").Value = -6799
$ws.Range("H122").Value = 3556.8572
$ws.Range("J122").Value = 4188.6665
$ws.Range("K122").Value = 9249
$ws.Range("L122").Value = 12565.9995
$ws.Range("N122").Value = -17465.9995
$ws.Range("I122").Value = 3083

# Row 34
$ws.Range("H34").Value = 108411.39
$ws.Range("M34").Value = -224263.22
$ws.Range("L34").Value = 22898.053
$ws.Range("I34").Value = 224465.22
$ws.Range("K34").Value = 224465.22
$ws.Range("N34").Value = -23302.053
$ws.Range("J34").Value = 22898.053

# Row 132
$ws.Range("J132").Value = 470728.8
$ws.Range("M132").Value = -142862150
$ws.Range("N132").Value = -1417246.4
$ws.Range("K132").Value = 142864680
$ws.Range("H132").Value = 31413462
$ws.Range("L132").Value = 1412186.4
$ws.Range("I132").Value = 47621560

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("L12").Value = 210
$ws.Range("H12").Value = 70
$ws.Range("N12").Value = -556
$ws.Range("J12").Value = 70

# Row 23
$ws.Range("I23").Value = 46.923077
$ws.Range("K23").Value = 140.769231
$ws.Range("M23").Value = 94.23076900000001
$ws.Range("H23").Value = 279.8421

$ws = $wb.Worksheets.Item("GSM")
# Row 136
$ws.Range("H136").Value = 25082.3
$ws.Range("N136").Value = -80346.89999999999
$ws.Range("J136").Value = 25082.3
$ws.Range("L136").Value = 75246.89999999999

# Row 46
$ws.Range("H46").Value = 16912.5
$ws.Range("M46").Value = -6169
$ws.Range("J46").Value = 27500
$ws.Range("N46").Value = -27812
$ws.Range("I46").Value = 6325
$ws.Range("K46").Value = 6325
$ws.Range("L46").Value = 27500

# Row 126
$ws.Range("L126").Value = 11414.0001
$ws.Range("J126").Value = 3804.6667
$ws.Range("H126").Value = 1193499.4
$ws.Range("N126").Value = -16354.0001

# Row 43
$ws.Range("M43").Value = -1827.1428
$ws.Range("K43").Value = 1978.1428
$ws.Range("I43").Value = 1978.1428
$ws.Range("H43").Value = 6837.9

# Row 122
$ws.Range("M122").Value = -204407.41
$ws.Range("H122").Value = 43666.617
$ws.Range("J122").Value = 9185.909
$ws.Range("K122").Value = 206857.41
$ws.Range("L122").Value = 27557.727
$ws.Range("N122").Value = -32457.727
$ws.Range("I122").Value = 68952.47

# Row 70
$ws.Range("I70").Value = 6549.909
$ws.Range("H70").Value = 6503.3335
$ws.Range("L70").Value = 5991
$ws.Range("N70").Value = -6531
$ws.Range("K70").Value = 6549.909
$ws.Range("M70").Value = -6279.909
$ws.Range("J70").Value = 5991

# Row 102
$ws.Range("H102").Value = 4021.5386
$ws.Range("I102").Value = 3083.1
$ws.Range("K102").Value = 3083.1
$ws.Range("M102").Value = -1461.1

# Row 73
$ws.Range("J73").Value = 5991
$ws.Range("K73").Value = 6549.909
$ws.Range("M73").Value = -5613.909
$ws.Range("N73").Value = -7863
$ws.Range("L73").Value = 5991
$ws.Range("I73").Value = 6549.909
$ws.Range("H73").Value = 6503.3335

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("I7").Value = 3903.6924
$ws.Range("K7").Value = 3903.6924
$ws.Range("M7").Value = -3791.6924
$ws.Range("J7").Value = 4300
$ws.Range("L7").Value = 4300
$ws.Range("N7").Value = -4524

# Row 16
$ws.Range("K16").Value = 687.6667
$ws.Range("H16").Value = 690.75
$ws.Range("I16").Value = 687.6667
$ws.Range("M16").Value = -517.6667

# Row 126
$ws.Range("L126").Value = 12900
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 11711.0772
$ws.Range("M126").Value = -9241.0772
$ws.Range("I126").Value = 3903.6924
$ws.Range("N126").Value = -17840

# Row 23
$ws.Range("L23").Value = 40000
$ws.Range("N23").Value = -40460
$ws.Range("J23").Value = 40000
$ws.Range("H23").Value = 33998.5

# Row 132
$ws.Range("J132").Value = 9128.429
$ws.Range("N132").Value = -32445.287
$ws.Range("H132").Value = 2907299
$ws.Range("L132").Value = 27385.287

# Row 29
$ws.Range("H29").Value = 6000
$ws.Range("N29").Value = -6590
$ws.Range("L29").Value = 6000
$ws.Range("J29").Value = 6000

# Row 40
$ws.Range("I40").Value = 5124.5
$ws.Range("H40").Value = 4756.857
$ws.Range("M40").Value = -4988.5
$ws.Range("K40").Value = 5124.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("M132").Value = -54890354
$ws.Range("K132").Value = 54892884
$ws.Range("H132").Value = 9156381
$ws.Range("I132").Value = 18297628

# Row 32
$ws.Range("I32").Value = 2500
$ws.Range("K32").Value = 2500
$ws.Range("H32").Value = 2500
$ws.Range("M32").Value = -2183

Write-Output "Applied all updates"